# Update "想去人数" (want-to-go count) figures in the F column across
# the 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets,
# matching the data refresh captured in commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibit.Range("F2").Value  = 15203
$wsExhibit.Range("F3").Value  = 19614
$wsExhibit.Range("F5").Value  = 182
$wsExhibit.Range("F14").Value = 228
$wsExhibit.Range("F15").Value = 258
$wsExhibit.Range("F16").Value = 78
$wsExhibit.Range("F17").Value = 1534
$wsExhibit.Range("F20").Value = 120
$wsExhibit.Range("F21").Value = 251
$wsExhibit.Range("F22").Value = 8278
$wsExhibit.Range("F26").Value = 72
$wsExhibit.Range("F27").Value = 1279
$wsExhibit.Range("F28").Value = 32
$wsExhibit.Range("F30").Value = 24
$wsExhibit.Range("F31").Value = 6603
$wsExhibit.Range("F32").Value = 139
$wsExhibit.Range("F33").Value = 80
$wsExhibit.Range("F36").Value = 318
$wsExhibit.Range("F37").Value = 5650
$wsExhibit.Range("F38").Value = 1021
$wsExhibit.Range("F39").Value = 34
$wsExhibit.Range("F41").Value = 67

# 演出 (sheet2)
$wsShow.Range("F3").Value = 29

# 全部类型 (sheet4)
$wsAll.Range("F2").Value  = 15203
$wsAll.Range("F3").Value  = 19614
$wsAll.Range("F5").Value  = 182
$wsAll.Range("F14").Value = 228
$wsAll.Range("F15").Value = 258
$wsAll.Range("F16").Value = 78
$wsAll.Range("F17").Value = 1534
$wsAll.Range("F21").Value = 120
$wsAll.Range("F22").Value = 251
$wsAll.Range("F23").Value = 8278
$wsAll.Range("F24").Value = 994
$wsAll.Range("F27").Value = 72
$wsAll.Range("F28").Value = 1279
$wsAll.Range("F29").Value = 32
$wsAll.Range("F31").Value = 24
$wsAll.Range("F32").Value = 29
$wsAll.Range("F34").Value = 6603
$wsAll.Range("F35").Value = 139
$wsAll.Range("F36").Value = 80
$wsAll.Range("F39").Value = 318
$wsAll.Range("F40").Value = 5650
$wsAll.Range("F41").Value = 1021
$wsAll.Range("F42").Value = 34
$wsAll.Range("F44").Value = 67
